$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4693.943
$ws.Range("J38").Value = 6550
$ws.Range("L38").Value = 19650
$ws.Range("N38").Value = -20394

$ws.Range("H107").Value = 2331.125
$ws.Range("J107").Value = 2346.4285
$ws.Range("L107").Value = 2346.4285
$ws.Range("N107").Value = -6186.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7383.857
$ws.Range("I2").Value = 7136.4
$ws.Range("J2").Value = 8002.5
$ws.Range("K2").Value = 7136.4
$ws.Range("L2").Value = 8002.5
$ws.Range("M2").Value = -7023.4
$ws.Range("N2").Value = -8228.5

$ws.Range("H32").Value = 4860.745
$ws.Range("I32").Value = 4547.875
$ws.Range("J32").Value = 9866.666999999999
$ws.Range("K32").Value = 4547.875
$ws.Range("L32").Value = 9866.666999999999
$ws.Range("M32").Value = -4260.875
$ws.Range("N32").Value = -10440.667

$ws.Range("H61").Value = 3201.8333
$ws.Range("I61").Value = 2905.0967
$ws.Range("J61").Value = 4038.0908
$ws.Range("K61").Value = 2905.0967
$ws.Range("L61").Value = 4038.0908
$ws.Range("M61").Value = -2693.0967
$ws.Range("N61").Value = -4462.0908

$ws.Range("H74").Value = 45552.5
$ws.Range("I74").Value = 48691.19
$ws.Range("J74").Value = 4749.5
$ws.Range("K74").Value = 48691.19
$ws.Range("L74").Value = 4749.5
$ws.Range("M74").Value = -47817.19
$ws.Range("N74").Value = -6497.5

$ws.Range("H77").Value = 45552.5
$ws.Range("I77").Value = 48691.19
$ws.Range("J77").Value = 4749.5
$ws.Range("K77").Value = 243455.95
$ws.Range("L77").Value = 23747.5
$ws.Range("M77").Value = -239087.95
$ws.Range("N77").Value = -32483.5

$ws.Range("H97").Value = 1377.125
$ws.Range("I97").Value = 1272.238
$ws.Range("J97").Value = 2111.3333
$ws.Range("K97").Value = 1272.238
$ws.Range("L97").Value = 2111.3333
$ws.Range("M97").Value = -776.2380000000001
$ws.Range("N97").Value = -3103.3333

$ws.Range("H102").Value = 3389.28
$ws.Range("I102").Value = 3034.9048
$ws.Range("K102").Value = 3034.9048
$ws.Range("M102").Value = -1412.9048

$ws.Range("H110").Value = 70817
$ws.Range("I110").Value = 70817
$ws.Range("K110").Value = 70817
$ws.Range("M110").Value = -68772

$ws.Range("H116").Value = 7383.857
$ws.Range("I116").Value = 7136.4
$ws.Range("J116").Value = 8002.5
$ws.Range("K116").Value = 7136.4
$ws.Range("L116").Value = 8002.5
$ws.Range("M116").Value = -4842.4
$ws.Range("N116").Value = -12590.5

$ws.Range("H132").Value = 23756.875
$ws.Range("I132").Value = 1822.1273
$ws.Range("K132").Value = 5466.3819
$ws.Range("M132").Value = -2936.3819

$ws.Range("H136").Value = 3201.8333
$ws.Range("I136").Value = 2905.0967
$ws.Range("J136").Value = 4038.0908
$ws.Range("K136").Value = 8715.2901
$ws.Range("L136").Value = 12114.2724
$ws.Range("M136").Value = -6165.2901
$ws.Range("N136").Value = -17214.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7383.857
$ws.Range("I3").Value = 7136.4
$ws.Range("J3").Value = 8002.5
$ws.Range("K3").Value = 7136.4
$ws.Range("L3").Value = 8002.5
$ws.Range("M3").Value = -7022.4
$ws.Range("N3").Value = -8230.5

$ws.Range("H20").Value = 1681.4
$ws.Range("J20").Value = 1719.6
$ws.Range("L20").Value = 1719.6
$ws.Range("N20").Value = -2213.6

$ws.Range("H86").Value = 23705.785
$ws.Range("I86").Value = 13186.111
$ws.Range("K86").Value = 13186.111
$ws.Range("M86").Value = -12063.111

$ws.Range("H89").Value = 23705.785
$ws.Range("I89").Value = 13186.111
$ws.Range("K89").Value = 65930.55500000001
$ws.Range("M89").Value = -60314.55500000001

$ws.Range("H105").Value = 23292.3
$ws.Range("I105").Value = 24769.223
$ws.Range("K105").Value = 24769.223
$ws.Range("M105").Value = -23022.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4800.125
$ws.Range("I31").Value = 2841.889
$ws.Range("J31").Value = 5975.067
$ws.Range("K31").Value = 2841.889
$ws.Range("L31").Value = 5975.067
$ws.Range("M31").Value = -2546.889
$ws.Range("N31").Value = -6565.067

$ws.Range("H34").Value = 4800.125
$ws.Range("I34").Value = 2841.889
$ws.Range("J34").Value = 5975.067
$ws.Range("K34").Value = 2841.889
$ws.Range("L34").Value = 5975.067
$ws.Range("M34").Value = -2639.889
$ws.Range("N34").Value = -6379.067

$ws.Range("H134").Value = 3946.7693
$ws.Range("I134").Value = 3811.7778
$ws.Range("K134").Value = 11435.3334
$ws.Range("M134").Value = -8900.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 7500
$ws.Range("J82").Value = 7500
$ws.Range("L82").Value = 22500
$ws.Range("N82").Value = -23312

$ws.Range("H85").Value = 7500
$ws.Range("J85").Value = 7500
$ws.Range("L85").Value = 22500
$ws.Range("N85").Value = -25308

$ws.Range("H131").Value = 1284.9166
$ws.Range("J131").Value = 1658.6154
$ws.Range("L131").Value = 4975.8462
$ws.Range("N131").Value = -15055.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5405.269
$ws.Range("I102").Value = 2596.6875
$ws.Range("K102").Value = 2596.6875
$ws.Range("M102").Value = -974.6875

$ws.Range("H113").Value = 2767.6
$ws.Range("I113").Value = 2043.5
$ws.Range("K113").Value = 2043.5
$ws.Range("M113").Value = 126.5

$ws.Range("H126").Value = 1429.6666
$ws.Range("I126").Value = 1495
$ws.Range("J126").Value = 1299
$ws.Range("K126").Value = 4485
$ws.Range("L126").Value = 3897
$ws.Range("M126").Value = -2015
$ws.Range("N126").Value = -8837

$ws.Range("H132").Value = 3330
$ws.Range("I132").Value = 3266.6667
$ws.Range("J132").Value = 3357.1428
$ws.Range("K132").Value = 9800.000100000001
$ws.Range("L132").Value = 10071.4284
$ws.Range("M132").Value = -7270.000100000001
$ws.Range("N132").Value = -15131.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3562.625
$ws.Range("I40").Value = 3666.6667
$ws.Range("J40").Value = 3500.2
$ws.Range("K40").Value = 3666.6667
$ws.Range("L40").Value = 3500.2
$ws.Range("M40").Value = -3530.6667
$ws.Range("N40").Value = -3772.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 169996
$ws.Range("J46").Value = 169996
$ws.Range("L46").Value = 169996
$ws.Range("N46").Value = -170458

$ws.Range("H126").Value = 4655.524
$ws.Range("I126").Value = 4795.533
$ws.Range("J126").Value = 4305.5
$ws.Range("K126").Value = 14386.599
$ws.Range("L126").Value = 12916.5
$ws.Range("M126").Value = -11916.599
$ws.Range("N126").Value = -17856.5

$ws.Range("H134").Value = 169996
$ws.Range("J134").Value = 169996
$ws.Range("L134").Value = 509988
$ws.Range("N134").Value = -515058
